$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to remain a Text cell (matches the source workbook,
    # which stores every one of these values as an inline string) even when
    # the text looks like a pure number (COM normally auto-converts those).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range("D2").Value = "63.417.50"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "3.068.64"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  -0.31%  "
Set-TextValue $ws.Range("D5") "589.23"
$ws.Range("E5").Value = "  -0.56%  "
Set-TextValue $ws.Range("D6") "155.57"
$ws.Range("E6").Value = "  +4.99%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "3.068.45"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("E10").Value = "  -3.81%  "
Set-TextValue $ws.Range("D11") "5.83"
$ws.Range("E11").Value = "  -1.38%  "
Set-TextValue $ws.Range("D12") "0.449"
$ws.Range("E12").Value = "  -2.14%  "
Set-TextValue $ws.Range("D13") "36.97"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("E14").Value = "  -4.28%  "
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").Value = "3.573.32"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "63.459.59"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").Value = "3.066.01"
$ws.Range("E19").Value = "  -2.74%  "
Set-TextValue $ws.Range("D20") "472.78"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  -1.64%  "
Set-TextValue $ws.Range("D22") "0.704"
$ws.Range("E22").Value = "  -4.37%  "
Set-TextValue $ws.Range("D23") "7.51"
$ws.Range("E23").Value = "  -1.86%  "
Set-TextValue $ws.Range("D24") "2.41"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D25") "12.82"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D26") "80.64"
$ws.Range("E26").Value = "  -0.97%  "
Set-TextValue $ws.Range("D27") "10.35"
$ws.Range("E27").Value = "  +2.43%  "
Set-TextValue $ws.Range("D28") "0.998"
$ws.Range("E28").Value = "  -0.21%  "
Set-TextValue $ws.Range("D29") "7.42"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -2.55%  "
Set-TextValue $ws.Range("D32") "2.15"
$ws.Range("E32").Value = "  -3.56%  "
Set-TextValue $ws.Range("D33") "0.111"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("D35").Value = "0.0₃0818"
$ws.Range("E35").Value = "  -5.04%  "
$ws.Range("E36").Value = "  -2.13%  "
Set-TextValue $ws.Range("D37") "3.30"
$ws.Range("E37").Value = "  +0.21%  "
Set-TextValue $ws.Range("D38") "5.98"
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("E39").Value = "  -4.70%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D40") "50.51"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D41") "9.19"
$ws.Range("E41").Value = "  -1.26%  "
Set-TextValue $ws.Range("D42") "436.55"
$ws.Range("E42").Value = "  -6.31%  "
Set-TextValue $ws.Range("D43") "0.285"
$ws.Range("E43").Value = "  -3.60%  "
Set-TextValue $ws.Range("D44") "40.76"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("E46").Value = "  -4.32%  "
$ws.Range("D47").Value = "2.793.77"
$ws.Range("E47").Value = "  -3.72%  "
Set-TextValue $ws.Range("D48") "129.69"
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("E49").Value = "  +0.00%  "
Set-TextValue $ws.Range("D50") "25.13"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("E51").Value = "  -0.80%  "
